$d = $word.ActiveDocument

# 1. Re-id the "smarthosting" bookmark from 1 to 0 by deleting and re-adding it
#    (the engine assigns the lowest free w:id, which becomes 0 once the doc's
#    other bookmark id no longer exists / was never present here).
$bm = $d.Bookmarks("smarthosting")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("smarthosting", $bmRange)

# 2. Update the lead-in sentence before the exchange links and drop the
#    "CryptoBridge" HYPERLINK field entirely (code + displayed text).
$d.Content.Find.Execute(
    "SmartCash можно купить на таких биржах, как, например, ", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Smart can be obtained from exchanges", 2)
$d.Fields.Item(1).Delete()

# 3. Turn the " или " connector into a lone "." and drop the "HitBTC"
#    HYPERLINK field entirely (code + displayed text).
$d.Content.Find.Execute(
    " или ", $true, $false, $false, $false, $false, $true, 1, $false,
    ".", 2)
$d.Fields.Item(1).Delete()

# 4. Remove the now-orphaned bold "." run left over from the old sentence,
#    and update the remaining lead-in text for the "full list" hyperlink.
$d.Content.Find.Execute(
    ".Полный список бирж, где можно купить SmartCash,  ", $true, $false,
    $false, $false, $false, $true, 1, $false,
    " For the full list of places to obtain Smart click ", 2)

# 5. Rename the "здесь" hyperlink display text to "here".
$h = $d.Hyperlinks.Item(1)
$h.TextToDisplay = "here"
